# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-05-28 Wednesday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-05-29 Thursday", 2)

# Update the arithmetic-problem table cells by explicit (row, column)
# position, since several of the new values collide with old values used
# elsewhere in the table (a naive global find/replace would cascade).
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "20÷3=6, 2"
$t.Cell(1,2).Range.Text = "11÷3=3, 2"
$t.Cell(1,3).Range.Text = "59÷2=29, 1"
$t.Cell(1,4).Range.Text = "19÷6=3, 1"
$t.Cell(1,5).Range.Text = "77÷2=38, 1"

$t.Cell(5,1).Range.Text = "84÷3=28, 0"
$t.Cell(5,2).Range.Text = "74÷2=37, 0"
$t.Cell(5,3).Range.Text = "46÷2=23, 0"
$t.Cell(5,4).Range.Text = "11÷7=1, 4"
$t.Cell(5,5).Range.Text = "42÷8=5, 2"

$t.Cell(9,1).Range.Text = "81÷3=27, 0"
$t.Cell(9,2).Range.Text = "23÷3=7, 2"
$t.Cell(9,3).Range.Text = "96÷6=16, 0"
$t.Cell(9,4).Range.Text = "26÷9=2, 8"
$t.Cell(9,5).Range.Text = "10÷3=3, 1"

$t.Cell(13,1).Range.Text = "72÷6=12, 0"
$t.Cell(13,2).Range.Text = "31÷3=10, 1"
$t.Cell(13,3).Range.Text = "82÷9=9, 1"
$t.Cell(13,4).Range.Text = "57÷8=7, 1"
$t.Cell(13,5).Range.Text = "43÷8=5, 3"

$t.Cell(17,1).Range.Text = "87÷3=29, 0"
$t.Cell(17,2).Range.Text = "23÷4=5, 3"
$t.Cell(17,3).Range.Text = "43÷9=4, 7"
$t.Cell(17,4).Range.Text = "89÷2=44, 1"
$t.Cell(17,5).Range.Text = "71÷9=7, 8"

Write-Host "applied edits"
